# "Add files via upload" — re-upload of "R in Jupyter.pptx" whose only
# user-visible content change is the corrected GitHub repo URL/text shown
# on the "download demo notebook" slide (slide 12, the text-box shape that
# also carries the http://blog.revolutionanalytics.com... link).
#
# Old repo name: DaveSnell/Demo-of-R-in-Jupyter-notebook
# New repo name: DaveSnell/demo-of-R-in-Jupyter
#
# We rewrite only the run text that contains the old URL fragment, using
# TextRange.Characters(start, length) so the surrounding runs (and the
# <a:hlinkClick> on this very run) are left completely untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

$oldText = "github.com/DaveSnell/Demo-of-R-in-Jupyter-notebook"
$newText = "github.com/DaveSnell/demo-of-R-in-Jupyter"

$full = $tr.Text
$idx = $full.IndexOf($oldText)
if ($idx -ge 0) {
    $target = $tr.Characters($idx + 1, $oldText.Length)
    $target.Text = $newText
}
